$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2 through 12
# from serial date 45175 (2023-09-06) to 45183 (2023-09-14)
foreach ($row in 2..12) {
    $ws.Cells.Item($row, 3).Value = 45183
}
